$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B279").Value = 183.3866

$dates = @("06-10-2021","07-10-2021","08-10-2021","09-10-2021","10-10-2021","11-10-2021","12-10-2021")
$values = @(186.0993, 186.4168, 188.7231, 187.4273, 185.2072, 188.079, 197.1793)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 280 + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
